# Re-sort the classification table by the new average lap time (Lap_Time),
# recomputed after fixing the number_of_laps x average_lap_time chart.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 'Virtualdo Pereira'
$ws.Cells.Item(2, 2).Value = 'Brazil'
$ws.Cells.Item(2, 3).Value = 'ACS Racing'
$ws.Cells.Item(2, 4).Value = 13
$ws.Cells.Item(2, 5).Value = 'Medium'
$ws.Cells.Item(2, 6).Value = '0:01:07.600338'

$ws.Cells.Item(3, 1).Value = 'George Nicholson'
$ws.Cells.Item(3, 2).Value = 'England'
$ws.Cells.Item(3, 3).Value = 'Empire GP'
$ws.Cells.Item(3, 4).Value = 4
$ws.Cells.Item(3, 5).Value = 'Medium'
$ws.Cells.Item(3, 6).Value = '0:01:07.888590'

$ws.Cells.Item(4, 1).Value = 'Giorgio Bussagna'
$ws.Cells.Item(4, 2).Value = 'Italy'
$ws.Cells.Item(4, 3).Value = 'Empire GP'
$ws.Cells.Item(4, 4).Value = 14
$ws.Cells.Item(4, 5).Value = 'Medium'
$ws.Cells.Item(4, 6).Value = '0:01:07.976509'

$ws.Cells.Item(5, 1).Value = 'Umineko Portela'
$ws.Cells.Item(5, 2).Value = 'Brazil'
$ws.Cells.Item(5, 3).Value = 'CoperSucca'
$ws.Cells.Item(5, 4).Value = 66
$ws.Cells.Item(5, 5).Value = 'Medium'
$ws.Cells.Item(5, 6).Value = '0:01:08.038231'

$ws.Cells.Item(6, 1).Value = 'Abílio de Souza'
$ws.Cells.Item(6, 2).Value = 'Equatorial Guinea'
$ws.Cells.Item(6, 3).Value = 'Gurgel RP'
$ws.Cells.Item(6, 4).Value = 16
$ws.Cells.Item(6, 5).Value = 'Medium'
$ws.Cells.Item(6, 6).Value = '0:01:08.223204'

$ws.Cells.Item(7, 1).Value = 'Leonardo Henrique'
$ws.Cells.Item(7, 2).Value = 'Brazil'
$ws.Cells.Item(7, 3).Value = 'CoperSucca'
$ws.Cells.Item(7, 4).Value = 24
$ws.Cells.Item(7, 5).Value = 'Medium'
$ws.Cells.Item(7, 6).Value = '0:01:08.254064'

$ws.Cells.Item(8, 1).Value = 'Chic Kane'
$ws.Cells.Item(8, 2).Value = 'England'
$ws.Cells.Item(8, 3).Value = 'Dani Ela Racing'
$ws.Cells.Item(8, 4).Value = 42
$ws.Cells.Item(8, 5).Value = 'Medium'
$ws.Cells.Item(8, 6).Value = '0:01:08.287349'

$ws.Cells.Item(9, 1).Value = 'Chu Pak-UI'
$ws.Cells.Item(9, 2).Value = 'North-Korea'
$ws.Cells.Item(9, 3).Value = 'Missuga Motors'
$ws.Cells.Item(9, 4).Value = 69
$ws.Cells.Item(9, 5).Value = 'Medium'
$ws.Cells.Item(9, 6).Value = '0:01:08.383341'

$ws.Cells.Item(10, 1).Value = 'Öster Tasion'
$ws.Cells.Item(10, 2).Value = 'Austria'
$ws.Cells.Item(10, 3).Value = 'Dani Ela Racing'
$ws.Cells.Item(10, 4).Value = 10
$ws.Cells.Item(10, 5).Value = 'Medium'
$ws.Cells.Item(10, 6).Value = '0:01:08.503860'

$ws.Cells.Item(11, 1).Value = 'Jorge Pelado'
$ws.Cells.Item(11, 2).Value = 'Equatorial Guinea'
$ws.Cells.Item(11, 3).Value = 'Missuga Motors'
$ws.Cells.Item(11, 4).Value = 65
$ws.Cells.Item(11, 5).Value = 'Medium'
$ws.Cells.Item(11, 6).Value = '0:01:08.520360'

$ws.Cells.Item(12, 1).Value = 'Max Overseas'
$ws.Cells.Item(12, 2).Value = 'England'
$ws.Cells.Item(12, 3).Value = 'Scuderia Archi Del''lappa'
$ws.Cells.Item(12, 4).Value = 71
$ws.Cells.Item(12, 5).Value = 'Hard'
$ws.Cells.Item(12, 6).Value = '0:01:08.523868'

$ws.Cells.Item(13, 1).Value = 'Thravekis Galludis'
$ws.Cells.Item(13, 2).Value = 'Greece'
$ws.Cells.Item(13, 3).Value = 'Gurgel RP'
$ws.Cells.Item(13, 4).Value = 77
$ws.Cells.Item(13, 5).Value = 'Medium'
$ws.Cells.Item(13, 6).Value = '0:01:08.562919'

$ws.Cells.Item(14, 1).Value = 'Marcelo Mastroianni'
$ws.Cells.Item(14, 2).Value = 'Italiano'
$ws.Cells.Item(14, 3).Value = 'ACS Racing'
$ws.Cells.Item(14, 4).Value = 84
$ws.Cells.Item(14, 5).Value = 'Medium'
$ws.Cells.Item(14, 6).Value = '0:01:08.827894'

$ws.Cells.Item(15, 1).Value = 'Chavez Tigrón'
$ws.Cells.Item(15, 2).Value = 'Mexico'
$ws.Cells.Item(15, 3).Value = 'Scuderia Archi Del''lappa'
$ws.Cells.Item(15, 4).Value = 19
$ws.Cells.Item(15, 5).Value = 'Hard'
$ws.Cells.Item(15, 6).Value = '0:01:08.914482'

$ws.Cells.Item(16, 1).Value = 'Kahn Templani Efdoux'
$ws.Cells.Item(16, 2).Value = 'Belgium'
$ws.Cells.Item(16, 3).Value = 'Aoi Yu Racers'
$ws.Cells.Item(16, 4).Value = 68
$ws.Cells.Item(16, 5).Value = 'Hard'
$ws.Cells.Item(16, 6).Value = '0:01:09.006399'

$ws.Cells.Item(17, 1).Value = 'Acistino Effoum'
$ws.Cells.Item(17, 2).Value = 'Belgium'
$ws.Cells.Item(17, 3).Value = 'Aoi Yu Racers'
$ws.Cells.Item(17, 4).Value = 23
$ws.Cells.Item(17, 5).Value = 'Hard'
$ws.Cells.Item(17, 6).Value = '0:01:09.258644'
